$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 636, pushing the existing rows 636-727
# down to 638-729 (this is how the new weekly records were added upstream,
# shifting the whole historical block down by two rows).
$ws.Rows.Item(636).Insert()
$ws.Rows.Item(636).Insert()

# Populate the two newly inserted rows (636 and 637) with the new records.
# Columns A,B,C,E,F,G,H,N,Q,R carry the same constant values as every other
# row in this data block.
$ws.Cells.Item(636, 1).Value = 6
$ws.Cells.Item(636, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(636, 3).Value = "Metropolitana"
$ws.Cells.Item(636, 4).Value = 45127
$ws.Cells.Item(636, 5).Value = 13
$ws.Cells.Item(636, 6).Value = 100112039
$ws.Cells.Item(636, 7).Value = "Ciboulette"
$ws.Cells.Item(636, 8).Value = "Sin especificar"
$ws.Cells.Item(636, 9).Value = "Primera"
$ws.Cells.Item(636, 10).Value = 650
$ws.Cells.Item(636, 11).Value = 1000
$ws.Cells.Item(636, 12).Value = 1200
$ws.Cells.Item(636, 13).Value = 1108
$ws.Cells.Item(636, 14).Value = "`$/docena de atados"
$ws.Cells.Item(636, 15).Value = "Región Metropolitana"
$ws.Cells.Item(636, 16).Value = 369
$ws.Cells.Item(636, 17).Value = 3
$ws.Cells.Item(636, 18).Value = "Hortaliza"

$ws.Cells.Item(637, 1).Value = 6
$ws.Cells.Item(637, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(637, 3).Value = "Metropolitana"
$ws.Cells.Item(637, 4).Value = 45127
$ws.Cells.Item(637, 5).Value = 13
$ws.Cells.Item(637, 6).Value = 100112039
$ws.Cells.Item(637, 7).Value = "Ciboulette"
$ws.Cells.Item(637, 8).Value = "Sin especificar"
$ws.Cells.Item(637, 9).Value = "Segunda"
$ws.Cells.Item(637, 10).Value = 220
$ws.Cells.Item(637, 11).Value = 800
$ws.Cells.Item(637, 12).Value = 800
$ws.Cells.Item(637, 13).Value = 800
$ws.Cells.Item(637, 14).Value = "`$/docena de atados"
$ws.Cells.Item(637, 15).Value = "Región Metropolitana"
$ws.Cells.Item(637, 16).Value = 267
$ws.Cells.Item(637, 17).Value = 3
$ws.Cells.Item(637, 18).Value = "Hortaliza"

# Make sure column D on the new rows keeps the date number format used by
# the rest of the column (style index 2 in the original file).
$ws.Cells.Item(636, 4).NumberFormat = $ws.Cells.Item(638, 4).NumberFormat
$ws.Cells.Item(637, 4).NumberFormat = $ws.Cells.Item(638, 4).NumberFormat
